# Auto-generated edit script applying the VD (column X) and CH (column AB) revisions
# described by the commit diff, across the Cases, Fatalities, Hospitalized and ICU sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Cases")
$ws.Range("X16").Value = 130
$ws.Range("AB16").Value = 533
$ws.Range("X17").Value = 203
$ws.Range("AB17").Value = 702
$ws.Range("X18").Value = 277
$ws.Range("AB18").Value = 1054
$ws.Range("X19").Value = 369
$ws.Range("AB19").Value = 1406
$ws.Range("X21").Value = 567
$ws.Range("AB21").Value = 2022
$ws.Range("X44").Value = 4494
$ws.Range("AB44").Value = 23143
$ws.Range("X45").Value = 4585
$ws.Range("AB45").Value = 23813
$ws.Range("X46").Value = 4682
$ws.Range("AB46").Value = 24481
$ws.Range("X47").Value = 4731
$ws.Range("AB47").Value = 24930
$ws.Range("X48").Value = 4771
$ws.Range("AB48").Value = 25393
$ws.Range("X49").Value = 4808
$ws.Range("AB49").Value = 25665
$ws.Range("X50").Value = 4852
$ws.Range("AB50").Value = 25916
$ws.Range("X51").Value = 4919
$ws.Range("AB51").Value = 26238
$ws.Range("X52").Value = 4968
$ws.Range("AB52").Value = 26556
$ws.Range("X53").Value = 4997
$ws.Range("AB53").Value = 26855
$ws.Range("X54").Value = 5029
$ws.Range("AB54").Value = 27167
$ws.Range("X55").Value = 5051
$ws.Range("AB55").Value = 27460
$ws.Range("X56").Value = 5062
$ws.Range("AB56").Value = 27637
$ws.Range("X57").Value = 5094
$ws.Range("AB57").Value = 27839
$ws.Range("X58").Value = 5126
$ws.Range("AB58").Value = 27998
$ws.Range("X59").Value = 5147
$ws.Range("AB59").Value = 28195
$ws.Range("X61").Value = 5202
$ws.Range("AB61").Value = 28601
$ws.Range("X62").Value = 5215
$ws.Range("AB62").Value = 28740
$ws.Range("X63").Value = 5226
$ws.Range("AB63").Value = 28820

$ws = $wb.Worksheets.Item("Fatalities")
$ws.Range("X54").Value = 303
$ws.Range("AB54").Value = 1398
$ws.Range("X55").Value = 311
$ws.Range("AB55").Value = 1440
$ws.Range("X56").Value = 316
$ws.Range("AB56").Value = 1464
$ws.Range("X57").Value = 327
$ws.Range("AB57").Value = 1505
$ws.Range("X58").Value = 335
$ws.Range("AB58").Value = 1546
$ws.Range("X59").Value = 339
$ws.Range("AB59").Value = 1572
$ws.Range("X60").Value = 348
$ws.Range("AB60").Value = 1601
$ws.Range("X61").Value = 350
$ws.Range("AB61").Value = 1629
$ws.Range("X62").Value = 360
$ws.Range("AB62").Value = 1653
$ws.Range("X63").Value = 362
$ws.Range("AB63").Value = 1667

$ws = $wb.Worksheets.Item("Hospitalized")
$ws.Range("X31").Value = 313
$ws.Range("AB31").Value = 1464
$ws.Range("X32").Value = 324
$ws.Range("AB32").Value = 1610
$ws.Range("X33").Value = 359
$ws.Range("AB33").Value = 1792
$ws.Range("X34").Value = 368
$ws.Range("AB34").Value = 1880
$ws.Range("X35").Value = 370
$ws.Range("AB35").Value = 1994
$ws.Range("X36").Value = 373
$ws.Range("AB36").Value = 2179
$ws.Range("X37").Value = 381
$ws.Range("AB37").Value = 2207
$ws.Range("X38").Value = 387
$ws.Range("AB38").Value = 2287
$ws.Range("X39").Value = 378
$ws.Range("AB39").Value = 2351
$ws.Range("X40").Value = 369
$ws.Range("AB40").Value = 2332
$ws.Range("X41").Value = 369
$ws.Range("AB41").Value = 2320
$ws.Range("X42").Value = 380
$ws.Range("AB42").Value = 2305
$ws.Range("X43").Value = 364
$ws.Range("AB43").Value = 2309
$ws.Range("X44").Value = 337
$ws.Range("AB44").Value = 2227
$ws.Range("X45").Value = 330
$ws.Range("AB45").Value = 2139
$ws.Range("X46").Value = 313
$ws.Range("AB46").Value = 2071
$ws.Range("X47").Value = 304
$ws.Range("AB47").Value = 2011
$ws.Range("X48").Value = 290
$ws.Range("AB48").Value = 1936
$ws.Range("X49").Value = 299
$ws.Range("AB49").Value = 1913
$ws.Range("X50").Value = 299
$ws.Range("AB50").Value = 1896
$ws.Range("X51").Value = 283
$ws.Range("AB51").Value = 1857
$ws.Range("X52").Value = 256
$ws.Range("AB52").Value = 1733
$ws.Range("X53").Value = 250
$ws.Range("AB53").Value = 1677
$ws.Range("X54").Value = 234
$ws.Range("AB54").Value = 1579
$ws.Range("X55").Value = 231
$ws.Range("AB55").Value = 1528
$ws.Range("X56").Value = 234
$ws.Range("AB56").Value = 1522
$ws.Range("X57").Value = 218
$ws.Range("AB57").Value = 1499
$ws.Range("X58").Value = 205
$ws.Range("AB58").Value = 1421
$ws.Range("X59").Value = 192
$ws.Range("AB59").Value = 1356
$ws.Range("X60").Value = 194
$ws.Range("AB60").Value = 1305
$ws.Range("X62").Value = 187
$ws.Range("AB62").Value = 1232
$ws.Range("X63").Value = 196
$ws.Range("AB63").Value = 1210

$ws = $wb.Worksheets.Item("ICU")
$ws.Range("X57").Value = 48
$ws.Range("AB57").Value = 258
$ws.Range("X58").Value = 47
$ws.Range("AB58").Value = 245
$ws.Range("X61").Value = 40
$ws.Range("AB61").Value = 200
$ws.Range("X62").Value = 40
$ws.Range("AB62").Value = 195
$ws.Range("X63").Value = 41
$ws.Range("AB63").Value = 186
